$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '243.56'
$ws.Range("B4").Value = 'LEO'
$ws.Range("C4").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '3.616'
$ws.Range("E4").Value = '3LEOLEO'
$ws.Range("B5").Value = 'HuobiToken'
$ws.Range("C5").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '5.396'
$ws.Range("E5").Value = '4HuobiTokenHT'
$ws.Range("B6").Value = 'Cronos'
$ws.Range("C6").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.05915'
$ws.Range("E6").Value = '5CronosCRO'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.454'
$ws.Range("E7").Value = '6GateTokenGT'
$ws.Range("B8").Value = 'KuCoinToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '6.553'
$ws.Range("E8").Value = '7KuCoinTokenKCS'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8136'
$ws.Range("E9").Value = '8MXTokenMX'
$ws.Range("B10").Value = 'FTXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9150'
$ws.Range("E10").Value = '9FTXTokenFTT'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1419'
$ws.Range("E11").Value = '10WazirXWRX'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07434'
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'
$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03275'
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03068'
$ws.Range("E14").Value = '13BitrueCoinBTR'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09342'
$ws.Range("E15").Value = '14BitMartTokenBMX'
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.873'
$ws.Range("E16").Value = '15MCDexMCB'
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.001571'
$ws.Range("E17").Value = '16BitForexTokenBF'
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.04682'
$ws.Range("E18").Value = '17CoinExTokenCET'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.005937'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0009873'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.00008603'
$ws.Range("B23").Value = 'BTSEToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.148'
$ws.Range("E23").Value = '22BTSETokenBTSE'
$ws.Range("B24").Value = 'One'
$ws.Range("C24").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.01133'
$ws.Range("E24").Value = '23OneONE'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0002901'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03961'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006192'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1075'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003001'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008626'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005179'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000750'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.8693'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002279'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002101'
